# Auto-generated edit script: applies market-data refresh values
# as described by the commit diff, per-sheet, per-cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5811.591
$ws.Range("I100").Value = 4327
$ws.Range("J100").Value = 7593.1
$ws.Range("K100").Value = 4327
$ws.Range("L100").Value = 7593.1
$ws.Range("M100").Value = -3786
$ws.Range("N100").Value = -8675.1
$ws.Range("H132").Value = 17108.094
$ws.Range("I132").Value = 1304.1923
$ws.Range("J132").Value = 85591.664
$ws.Range("K132").Value = 3912.5769
$ws.Range("L132").Value = 256774.992
$ws.Range("M132").Value = -1382.5769
$ws.Range("N132").Value = -261834.992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7070.263
$ws.Range("I32").Value = 6446.5356
$ws.Range("K32").Value = 6446.5356
$ws.Range("M32").Value = -6159.5356
$ws.Range("H88").Value = 13890740
$ws.Range("I88").Value = 2350
$ws.Range("J88").Value = 20834936
$ws.Range("K88").Value = 2350
$ws.Range("L88").Value = 20834936
$ws.Range("M88").Value = -1944
$ws.Range("N88").Value = -20835748
$ws.Range("H91").Value = 13890740
$ws.Range("I91").Value = 2350
$ws.Range("J91").Value = 20834936
$ws.Range("K91").Value = 2350
$ws.Range("L91").Value = 20834936
$ws.Range("M91").Value = -946
$ws.Range("N91").Value = -20837744
$ws.Range("H132").Value = 2814.9167
$ws.Range("I132").Value = 1660.9474
$ws.Range("K132").Value = 4982.8422
$ws.Range("M132").Value = -2452.8422

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3042.0667
$ws.Range("I20").Value = 2892.875
$ws.Range("J20").Value = 3212.5715
$ws.Range("K20").Value = 2892.875
$ws.Range("L20").Value = 3212.5715
$ws.Range("M20").Value = -2645.875
$ws.Range("N20").Value = -3706.5715
$ws.Range("H64").Value = 1528.4
$ws.Range("I64").Value = 1444
$ws.Range("J64").Value = 1549.5
$ws.Range("K64").Value = 1444
$ws.Range("L64").Value = 1549.5
$ws.Range("M64").Value = -1219
$ws.Range("N64").Value = -1999.5
$ws.Range("H67").Value = 1528.4
$ws.Range("I67").Value = 1444
$ws.Range("J67").Value = 1549.5
$ws.Range("K67").Value = 1444
$ws.Range("L67").Value = 1549.5
$ws.Range("M67").Value = -664
$ws.Range("N67").Value = -3109.5
$ws.Range("H86").Value = 3236.611
$ws.Range("I86").Value = 1868.7778
$ws.Range("J86").Value = 4604.4443
$ws.Range("K86").Value = 1868.7778
$ws.Range("L86").Value = 4604.4443
$ws.Range("M86").Value = -745.7778000000001
$ws.Range("N86").Value = -6850.4443
$ws.Range("H89").Value = 3236.611
$ws.Range("I89").Value = 1868.7778
$ws.Range("J89").Value = 4604.4443
$ws.Range("K89").Value = 9343.889000000001
$ws.Range("L89").Value = 23022.2215
$ws.Range("M89").Value = -3727.889000000001
$ws.Range("N89").Value = -34254.2215
$ws.Range("H94").Value = 5378370.5
$ws.Range("I94").Value = 1503
$ws.Range("K94").Value = 1503
$ws.Range("M94").Value = -1052
$ws.Range("H99").Value = 3330.5454
$ws.Range("I99").Value = 3150.625
$ws.Range("K99").Value = 3150.625
$ws.Range("M99").Value = -1652.625
$ws.Range("H105").Value = 3665.9312
$ws.Range("I105").Value = 2998.15
$ws.Range("J105").Value = 5149.8887
$ws.Range("K105").Value = 2998.15
$ws.Range("L105").Value = 5149.8887
$ws.Range("M105").Value = -1251.15
$ws.Range("N105").Value = -8643.8887
$ws.Range("H128").Value = 5519.5
$ws.Range("I128").Value = 5519.5
$ws.Range("K128").Value = 16558.5
$ws.Range("M128").Value = -14068.5
$ws.Range("H134").Value = 2327.6724
$ws.Range("I134").Value = 1530.6531
$ws.Range("K134").Value = 4591.9593
$ws.Range("M134").Value = -2056.9593

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3997.087
$ws.Range("I31").Value = 1936.9
$ws.Range("J31").Value = 5581.846
$ws.Range("K31").Value = 1936.9
$ws.Range("L31").Value = 5581.846
$ws.Range("M31").Value = -1641.9
$ws.Range("N31").Value = -6171.846
$ws.Range("H34").Value = 3997.087
$ws.Range("I34").Value = 1936.9
$ws.Range("J34").Value = 5581.846
$ws.Range("K34").Value = 1936.9
$ws.Range("L34").Value = 5581.846
$ws.Range("M34").Value = -1734.9
$ws.Range("N34").Value = -5985.846
$ws.Range("H86").Value = 20947.74
$ws.Range("I86").Value = 23426.438
$ws.Range("J86").Value = 17342.363
$ws.Range("K86").Value = 23426.438
$ws.Range("L86").Value = 17342.363
$ws.Range("M86").Value = -22303.438
$ws.Range("N86").Value = -19588.363
$ws.Range("H89").Value = 20947.74
$ws.Range("I89").Value = 23426.438
$ws.Range("J89").Value = 17342.363
$ws.Range("K89").Value = 117132.19
$ws.Range("L89").Value = 86711.815
$ws.Range("M89").Value = -111516.19
$ws.Range("N89").Value = -97943.815
$ws.Range("H107").Value = 11173.7
$ws.Range("I107").Value = 1253
$ws.Range("K107").Value = 1253
$ws.Range("M107").Value = 667
$ws.Range("H132").Value = 2276.4736
$ws.Range("I132").Value = 1912.0883
$ws.Range("K132").Value = 5736.2649
$ws.Range("M132").Value = -3206.2649
$ws.Range("H134").Value = 4290
$ws.Range("I134").Value = 3076.4866
$ws.Range("K134").Value = 9229.459800000001
$ws.Range("M134").Value = -6694.459800000001
$ws.Range("H141").Value = 145057.94
$ws.Range("J141").Value = 145057.94
$ws.Range("L141").Value = 145057.94
$ws.Range("N141").Value = -155417.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 92324.07000000001
$ws.Range("I70").Value = 104378.664
$ws.Range("K70").Value = 104378.664
$ws.Range("M70").Value = -104108.664
$ws.Range("H73").Value = 92324.07000000001
$ws.Range("I73").Value = 104378.664
$ws.Range("K73").Value = 104378.664
$ws.Range("M73").Value = -103442.664
$ws.Range("H122").Value = 3353.2
$ws.Range("I122").Value = 1288.7333
$ws.Range("J122").Value = 6449.9
$ws.Range("K122").Value = 3866.199900000001
$ws.Range("L122").Value = 19349.7
$ws.Range("M122").Value = -1416.199900000001
$ws.Range("N122").Value = -24249.7
$ws.Range("H132").Value = 2124.1
$ws.Range("I132").Value = 1972.7368
$ws.Range("K132").Value = 5918.2104
$ws.Range("M132").Value = -3388.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8051.107
$ws.Range("I40").Value = 9749.214
$ws.Range("J40").Value = 6353
$ws.Range("K40").Value = 9749.214
$ws.Range("L40").Value = 6353
$ws.Range("M40").Value = -9613.214
$ws.Range("N40").Value = -6625
$ws.Range("H93").Value = 2353.5
$ws.Range("I93").Value = 2141.3333
$ws.Range("J93").Value = 2626.2856
$ws.Range("K93").Value = 2141.3333
$ws.Range("L93").Value = 2626.2856
$ws.Range("M93").Value = -893.3332999999998
$ws.Range("N93").Value = -5122.2856
$ws.Range("H122").Value = 4701.0713
$ws.Range("I122").Value = 2488.75
$ws.Range("J122").Value = 7650.8335
$ws.Range("K122").Value = 7466.25
$ws.Range("L122").Value = 22952.5005
$ws.Range("M122").Value = -5016.25
$ws.Range("N122").Value = -27852.5005
$ws.Range("H132").Value = 4081.0417
$ws.Range("I132").Value = 2446.9092
$ws.Range("J132").Value = 5463.769
$ws.Range("K132").Value = 7340.7276
$ws.Range("L132").Value = 16391.307
$ws.Range("M132").Value = -4810.7276
$ws.Range("N132").Value = -21451.307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 33043.5
$ws.Range("I69").Value = 49000
$ws.Range("J69").Value = 29852.2
$ws.Range("K69").Value = 49000
$ws.Range("L69").Value = 29852.2
$ws.Range("M69").Value = -48251
$ws.Range("N69").Value = -31350.2
$ws.Range("H70").Value = 17500
$ws.Range("J70").Value = 17500
$ws.Range("L70").Value = 17500
$ws.Range("N70").Value = -18130
$ws.Range("H72").Value = 33043.5
$ws.Range("I72").Value = 49000
$ws.Range("J72").Value = 29852.2
$ws.Range("K72").Value = 147000
$ws.Range("L72").Value = 89556.60000000001
$ws.Range("M72").Value = -143256
$ws.Range("N72").Value = -97044.60000000001
$ws.Range("H73").Value = 17500
$ws.Range("J73").Value = 17500
$ws.Range("L73").Value = 17500
$ws.Range("N73").Value = -19684
$ws.Range("H81").Value = 15153711
$ws.Range("I81").Value = 2566.1667
$ws.Range("J81").Value = 33335086
$ws.Range("K81").Value = 5132.3334
$ws.Range("L81").Value = 66670172
$ws.Range("M81").Value = -4071.3334
$ws.Range("N81").Value = -66672294
$ws.Range("H84").Value = 15153711
$ws.Range("I84").Value = 2566.1667
$ws.Range("J84").Value = 33335086
$ws.Range("K84").Value = 25661.667
$ws.Range("L84").Value = 333350860
$ws.Range("M84").Value = -20357.667
$ws.Range("N84").Value = -333361468
$ws.Range("H126").Value = 2255
$ws.Range("I126").Value = 2161.875
$ws.Range("K126").Value = 6485.625
$ws.Range("M126").Value = -4015.625
